# Update "想去人数" (F column) values on the "展览" and "全部类型" sheets.
$wb = $excel.ActiveWorkbook

# Row -> new F-column value
$updates = @{
    3  = 2195
    4  = 84
    5  = 13047
    6  = 72
    10 = 1175
    12 = 13750
    13 = 14306
    21 = 32
    22 = 1088
    25 = 5370
    28 = 304
    29 = 12
    30 = 19
}

foreach ($sheetName in @("展览", "全部类型")) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Cells.Item($row, 6).Value = $updates[$row]
    }
}
